# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# This updates the "K" column (G2:G21) with newly computed strikeout counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 4
    4  = 4
    5  = 6
    6  = 4
    7  = 6
    8  = 6
    9  = 9
    10 = 6
    11 = 9
    12 = 10
    13 = 5
    14 = 2
    15 = 4
    16 = 3
    17 = 8
    18 = 5
    19 = 4
    20 = 2
    21 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
